$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.905.49"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.753.35"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").Value = "'335.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").Value = "'0.3401"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'46.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'1.115"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "'0.07226"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'6.165"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'7.147"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "1.755.69"
$ws.Range("D17").Value = "'0.00001060"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "'78.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D22").Value = "'6.227"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").Value = "27.909.68"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").Value = "'11.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("D25").Value = "'2.386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "'152.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'19.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("D28").Value = "'2.314"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.41%  "
$ws.Range("D29").Value = "1.956.68"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  -11.95%  "
$ws.Range("D31").Value = "'132.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "'4.024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'5.838"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("D34").Value = "'0.08804"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'12.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("D36").Value = "'0.6584"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").Value = "'0.02287"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.24%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.145"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06158"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'0.2104"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "'7.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "'13.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "'0.6075"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "'3.822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").Value = "'126.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("E49").Value = "  -5.16%  "
$ws.Range("D50").Value = "'1.172"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "'1.117"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.11%  "
